# Blind model names in the "Evaluations" sheet (column C, rows 2-49)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Evaluations")

$mapping = @{
    "claude-opus-4.5" = "Model A"
    "gemini-3-pro"    = "Model B"
    "gpt-5.1"         = "Model C"
    "kimi-k2"         = "Model D"
}

for ($row = 2; $row -le 49; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $current = $cell.Value()
    if ($mapping.ContainsKey($current)) {
        $cell.Value = $mapping[$current]
    }
}
